$d = $word.ActiveDocument

# --- Block J: originally paragraphs 26..27 ---
$rngStart = $d.Paragraphs(26).Range.Start
$rngEnd = $d.Paragraphs(27).Range.End
$rng = $d.Range($rngStart, $rngEnd)
$xml_J = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Anwesende: </w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>alle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Aufteilung der Sprechparts der Präsentation. Vorbereitung läuft gut.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>13. März: Aufnahme der Präsentation</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Anwesende:</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> alle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">14. März 2020: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Anwesende: </w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Dennis Grunenberg</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">, Natasha Dudler (Zuhause: </w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Melanie </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Svab</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Während Dennis und Natasha den ganzen Nachmittag hindurch einen Bug bezüglich der Verbindung mehrerer Clients ausfindig zu machen versuchten, schnitt Melanie die einzelnen Teile der gestrigen Präsentation zu einem coolen Filmchen zusammen.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:lastRenderedPageBreak/><w:t>15. März 2020: Konferenz über Skype von 13:00-16:00</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Anwesende: alle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Einzelne Arbeiten wurden aufgeteilt. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Melanie = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Netwerkprotokoll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Rohail</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> = Chat, Dennis = Verbindungsstücke der einzelnen </w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Bereiche</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>, Natasha = Einstellungen</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> des</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> Client</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>s (Namen, etc.)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Gemeinsam suchten wir nach einem Weg, wie Clients miteinander kommunizieren können und besprachen den weiteren Verlauf bis zum zweiten Meilenstein.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">16. März 2020: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Anwesend: Natasha Dudler</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Update unseres </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Timetables</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>, sowie Programmieren von Klassen für Informationen des Clients an den Server.</w:t></w:r></w:p>'
$rng.InsertXML($xml_J)

# --- Block I: originally paragraphs 24..24 ---
$rng = $d.Paragraphs(24).Range
$xml_I = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p>'
$rng.InsertXML($xml_I)

# --- Block H: originally paragraphs 22..23 ---
$rngStart = $d.Paragraphs(22).Range.Start
$rngEnd = $d.Paragraphs(23).Range.End
$rng = $d.Range($rngStart, $rngEnd)
$xml_H = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Anwesende: Natasha Dudler</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Nach etlichen Stunden des Termine-Verschiebens und des Verzweifelns steht</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>nun der Terminkalender fest. Jetzt müssen wir uns nur noch daran halten…</w:t></w:r></w:p>'
$rng.InsertXML($xml_H)

# --- Block G: originally paragraphs 21..21 ---
$rng = $d.Paragraphs(21).Range
$xml_G = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">5. März 2020: Kreieren des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Timetables</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> mit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>GanttProjekt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>, sowie des Tagebuches</w:t></w:r></w:p>'
$rng.InsertXML($xml_G)

# --- Block F: originally paragraphs 19..20 ---
$rngStart = $d.Paragraphs(19).Range.Start
$rngEnd = $d.Paragraphs(20).Range.End
$rng = $d.Range($rngStart, $rngEnd)
$xml_F = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Anwesende: Melanie </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Svab</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p>'
$rng.InsertXML($xml_F)

# --- Block E: originally paragraphs 18..18 ---
$rng = $d.Paragraphs(18).Range
$xml_E = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">4. März 2020: Erstellen </w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">von </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Mockups</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$rng.InsertXML($xml_E)

# --- Block D: originally paragraphs 13..13 ---
$rng = $d.Paragraphs(13).Range
$xml_D = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Anwesende: Dennis Grunenberg, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Rohail</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Gondal</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$rng.InsertXML($xml_D)

# --- Block C: originally paragraphs 9..9 ---
$rng = $d.Paragraphs(9).Range
$xml_C = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Jeder brachte jeweils drei Ideen mit, die wir ausführlich besprachen. Nun haben wir uns für ein Schildkrötenspiel namens «The Floor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>is</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> Java» entschieden, zu dem wir sämtliche Regeln und Details festlegten.</w:t></w:r></w:p>'
$rng.InsertXML($xml_C)

# --- Block B: originally paragraphs 7..7 ---
$rng = $d.Paragraphs(7).Range
$xml_B = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Anwesende</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> der Gruppe</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>alle</w:t></w:r></w:p>'
$rng.InsertXML($xml_B)

# --- Block A: originally paragraphs 1..3 ---
$rngStart = $d.Paragraphs(1).Range.Start
$rngEnd = $d.Paragraphs(3).Range.End
$rng = $d.Range($rngStart, $rngEnd)
$xml_A = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Tagebuch für das Spiel “The Floor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="de-CH"/></w:rPr><w:t>is</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> Java”</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">von Melanie </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Svab</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">, Dennis Grunenberg, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Rohail</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Gondal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> Natasha Dudler</w:t></w:r></w:p>'
$rng.InsertXML($xml_A)

Write-Host "Final paragraph count:" $d.Paragraphs.Count